# Common: Preparing support for cascaded object creation
# Adds Czech translation rows for the new "atomizer" and "vendor" lab
# create-dialog labels to the "Import" worksheet (sheet1.xml).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Import")

# Helper: stamp a new data row with the same look & feel (style) as the
# last existing data row (row 210) by copying its formatting first, then
# overwriting the three cell values. Positional params only -- this
# runtime's PowerShell does not bind named "-Param value" arguments.
function Add-TranslationRow($RowIndex, $Lang, $Label, $Translation) {
    $ws.Range("A210:C210").Copy()
    $ws.Range("A" + $RowIndex + ":C" + $RowIndex).PasteSpecial(-4122)

    $ws.Range("A" + $RowIndex).Value = $Lang
    $ws.Range("B" + $RowIndex).Value = $Label
    $ws.Range("C" + $RowIndex).Value = $Translation
}

# Rows 211-216 and 218 are written top to bottom; row 217
# ("lab.vendor.create.title") is written last on purpose so the new
# shared-string entries land in the exact same append order the original
# workbook uses (lab.vendor.create.subtitle's strings are first-used
# before lab.vendor.create.title's).
Add-TranslationRow 211 "cs" "lab.atomizer.tooltip.create" "Přidat atomizér"
Add-TranslationRow 212 "cs" "lab.atomizer.create.title" "Vytvořit atomizér"
Add-TranslationRow 213 "cs" "lab.atomizer.create.subtitle" "Přidejte nový atomizér, který tak bude přístupný ostatních a v buildech."
Add-TranslationRow 214 "cs" "lab.atomizer.name.label" "Jméno"
Add-TranslationRow 215 "cs" "lab.atomizer.vendorId.label" "Výrobce"
Add-TranslationRow 216 "cs" "lab.vendor.tooltip.create" "Založit výrobce"
Add-TranslationRow 218 "cs" "lab.vendor.create.subtitle" "Výrobci jsou dostupní přes celou aplikaci, např. u atomizérů, modů, drátů a dalšího."
Add-TranslationRow 217 "cs" "lab.vendor.create.title" "Založit výrobce"

# Reflect the updated view state: selection resting on the newly added
# first row's key cell (the view's scroll position/topLeftCell is window
# chrome that this host does not persist through the COM bridge, so it
# is intentionally left alone beyond the selection itself).
$ws.Activate() | Out-Null
$ws.Range("B211").Select() | Out-Null

Write-Output "Added 8 translation rows (211-218) to Import sheet"
